# Dispatch Priority by Elec Source.xlsx - "Copy in EU-2024-develop branch"
#
# 1) "About" sheet: cell A4 ("None") had a stray font style applied (s="6");
#    clear it back to the default/unstyled cell.
# 2) "DPbES" sheet: flip the dispatch-priority flags (B:AE) for four rows:
#      row 2  (hard coal)   0 -> 1
#      row 6  (hydro)       1 -> 0
#      row 10 (biomass)     1 -> 0
#      row 12 (petroleum)   0 -> 1
# 3) View state: DPbES becomes the active/selected sheet (its previous
#    topLeftCell scroll position is cleared) with B10:AE10 selected, while
#    the About sheet is no longer the tab-selected one.

$wb = $excel.ActiveWorkbook

# --- About sheet: remove the stray style from A4 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A4").ClearFormats()

# --- DPbES sheet: update the Boolean dispatch-priority values ---
$wsDpbes = $wb.Worksheets.Item("DPbES")
$wsDpbes.Range("B2:AE2").Value = 1
$wsDpbes.Range("B6:AE6").Value = 0
$wsDpbes.Range("B10:AE10").Value = 0
$wsDpbes.Range("B12:AE12").Value = 1

# --- View state: activate DPbES and select B10:AE10 ---
$null = $wsDpbes.Activate()
$null = $wsDpbes.Range("B10:AE10").Select()
